# Actualiza los horarios de la Linea 141 (scrape 06:04:29) en las 3 hojas
# del libro: LP1912, LP1912-215, 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2,1).Value = "Última actualización: 06:04:29"
$ws1.Cells.Item(3,1).Value = "Total filas: 31"

# Filas existentes cuya hora de scrap / minutos (y a veces horario) cambian.
$ws1.Cells.Item(14,1).Value = "06:04:29"
$ws1.Cells.Item(14,4).Value = 7

$ws1.Cells.Item(15,1).Value = "06:04:29"
$ws1.Cells.Item(15,4).Value = 10

$ws1.Cells.Item(16,1).Value = "06:04:29"
$ws1.Cells.Item(16,4).Value = 17

$ws1.Cells.Item(17,1).Value = "06:04:29"
$ws1.Cells.Item(17,4).Value = 23

$ws1.Cells.Item(18,1).Value = "06:04:29"
$ws1.Cells.Item(18,4).Value = 25

$ws1.Cells.Item(19,1).Value = "06:04:29"
$ws1.Cells.Item(19,4).Value = 27

$ws1.Cells.Item(20,1).Value = "06:04:29"
$ws1.Cells.Item(20,4).Value = 40

$ws1.Cells.Item(21,1).Value = "06:04:29"
$ws1.Cells.Item(21,4).Value = 42

$ws1.Cells.Item(22,1).Value = "06:04:29"
$ws1.Cells.Item(22,4).Value = 55

$ws1.Cells.Item(23,1).Value = "06:04:29"
$ws1.Cells.Item(23,2).Value = "07:04"
$ws1.Cells.Item(23,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(23,4).Value = 60

$ws1.Cells.Item(24,1).Value = "06:04:29"
$ws1.Cells.Item(24,2).Value = "07:05"
$ws1.Cells.Item(24,3).Value = "15_ABASTO"
$ws1.Cells.Item(24,4).Value = 61

$ws1.Cells.Item(25,1).Value = "06:04:29"
$ws1.Cells.Item(25,2).Value = "07:07"
$ws1.Cells.Item(25,3).Value = "225_GOMEZ"
$ws1.Cells.Item(25,4).Value = 63

$ws1.Cells.Item(26,1).Value = "06:04:29"
$ws1.Cells.Item(26,2).Value = "07:11"
$ws1.Cells.Item(26,3).Value = "215A_EL PATO"
$ws1.Cells.Item(26,4).Value = 67

$ws1.Cells.Item(27,1).Value = "06:04:29"
$ws1.Cells.Item(27,2).Value = "07:15"
$ws1.Cells.Item(27,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(27,4).Value = 71

# Filas nuevas al final de la hoja (28-36).
$sheet1NewRows = @(
    @(28, "06:04:29", "07:21", "26_HERNANDEZ", 77, "LP1912"),
    @(29, "06:04:29", "07:23", "10_OLMOS", 79, "LP1912"),
    @(30, "06:04:29", "07:31", "11_ETCHEVERRY", 87, "LP1912"),
    @(31, "06:04:29", "07:31", "16_SANTA ANA", 87, "LP1912"),
    @(32, "06:04:29", "07:32", "84_COLONIA URQUIZA-ESC 49", 88, "LP1912"),
    @(33, "06:04:29", "07:36", "27_EL RETIRO", 92, "LP1912"),
    @(34, "06:04:29", "07:39", "10_OLMOS", 95, "LP1912"),
    @(35, "06:04:29", "07:47", "14_ABASTO", 103, "LP1912"),
    @(36, "06:04:29", "07:51", "215D_EL PATO", 107, "LP1912")
)
foreach ($row in $sheet1NewRows) {
    $r = $row[0]
    $ws1.Cells.Item($r,1).Value = $row[1]
    $ws1.Cells.Item($r,2).Value = $row[2]
    $ws1.Cells.Item($r,3).Value = $row[3]
    $ws1.Cells.Item($r,4).Value = $row[4]
    $ws1.Cells.Item($r,5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2,1).Value = "Última actualización: 06:04:29"
$ws2.Cells.Item(3,1).Value = "Total filas: 6"

$ws2.Cells.Item(8,1).Value = "06:04:29"
$ws2.Cells.Item(8,4).Value = 7

$ws2.Cells.Item(9,1).Value = "06:04:29"
$ws2.Cells.Item(9,4).Value = 42

$ws2.Cells.Item(10,1).Value = "06:04:29"
$ws2.Cells.Item(10,4).Value = 67

$ws2.Cells.Item(11,1).Value = "06:04:29"
$ws2.Cells.Item(11,2).Value = "07:51"
$ws2.Cells.Item(11,3).Value = "215D_EL PATO"
$ws2.Cells.Item(11,4).Value = 107
$ws2.Cells.Item(11,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Hoja 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2,1).Value = "Última actualización: 06:04:29"
$ws3.Cells.Item(3,1).Value = "Total filas: 8"

$ws3.Cells.Item(9,1).Value = "06:04:29"
$ws3.Cells.Item(9,4).Value = 5

$ws3.Cells.Item(11,1).Value = "06:04:29"
$ws3.Cells.Item(11,4).Value = 29

$ws3.Cells.Item(12,1).Value = "06:04:29"
$ws3.Cells.Item(12,4).Value = 56

$ws3.Cells.Item(13,1).Value = "06:04:29"
$ws3.Cells.Item(13,2).Value = "07:35"
$ws3.Cells.Item(13,3).Value = "215A_LA PLATA"
$ws3.Cells.Item(13,4).Value = 91
$ws3.Cells.Item(13,5).Value = "L6173"
